$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.343.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.312.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.40'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '577.97'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.34%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.98%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.411'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.884.40'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.86%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.564.76'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.98%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.320.34'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '443.45'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.72'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.54'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.77'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.89'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.456.37'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.13%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.76%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.49%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.81%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.89'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.35'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.24'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.02%  '

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.81'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.68%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.95'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.24'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.50%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.787.36'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.36%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.790'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.40%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.25'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.90'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.73%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.56%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.09'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '329.00'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.53%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.30%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.22'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.52%  '
